$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The report previously listed 2 short-stock items (rows 7-8) plus a totals
# row (9) and a footer row (10). The new export instead lists 6 items
# (rows 7-12), pushing the totals row to 13 and the footer to 14, and bumps
# the generation timestamp shown in the footer.
# ---------------------------------------------------------------------------

# Insert 4 fresh rows above the old totals row (row 9) so the old row 9
# (totals) becomes row 13 and the old row 10 (footer) becomes row 14. The
# inserted rows inherit neighbouring formatting, which we'll overwrite below.
$ws.Range("A9:Q12").Insert(-4121) | Out-Null

# Copy the exact cell formatting (borders/fills/fonts/number formats) used
# by the two existing item rows onto the four new item rows: row 7's look
# (odd item) -> rows 9 & 11, row 8's look (even item) -> rows 10 & 12.
$ws.Range("A7:Q7").Copy() | Out-Null
$ws.Range("A9:Q9").PasteSpecial(-4122, -4142, $false, $false) | Out-Null
$ws.Range("A11:Q11").PasteSpecial(-4122, -4142, $false, $false) | Out-Null

$ws.Range("A8:Q8").Copy() | Out-Null
$ws.Range("A10:Q10").PasteSpecial(-4122, -4142, $false, $false) | Out-Null
$ws.Range("A12:Q12").PasteSpecial(-4122, -4142, $false, $false) | Out-Null

$excel.CutCopyMode = 0

# Row heights for the new item rows (matches the alternating pattern already
# used by rows 7/8).
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5

# Re-create the same merge layout as rows 7/8 for each new item row.
foreach ($r in 9..12) {
    $ws.Range("A$r`:B$r").Merge()
    $ws.Range("C$r`:G$r").Merge()
    $ws.Range("H$r`:K$r").Merge()
    $ws.Range("L$r`:M$r").Merge()
    $ws.Range("N$r`:O$r").Merge()
}

# Helper: write a value as literal text (no forced numeric/date coercion),
# without disturbing the cell's existing number format.
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $fmt = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = $fmt
}

# --- Row 8 (item 2) now reports a different drug/price than before -------
Set-TextValue "C8" "GLIMEPIRIDE PLUS 4/30 MG 30 SCORED TAB."
Set-TextValue "N8" "123.00"
Set-TextValue "P8" "123.0000"

# --- Row 9 (item 3, new) ---------------------------------------------------
$ws.Range("A9").Value = 3
Set-TextValue "C9" "PANADOL ADVANCE 500 MG 48 TABLETS"
Set-TextValue "H9" "2:1"
Set-TextValue "L9" "1"
Set-TextValue "N9" "92.00"
Set-TextValue "P9" "23.0000"
Set-TextValue "Q9" "0:1"

# --- Row 10 (item 4, new) --------------------------------------------------
$ws.Range("A10").Value = 4
Set-TextValue "C10" "SPASMOFREE 5MG/2ML I.V./I.M. 3 AMP."
Set-TextValue "H10" "2:1"
Set-TextValue "L10" "1"
Set-TextValue "N10" "54.00"
Set-TextValue "P10" "17.8200"
Set-TextValue "Q10" "0:1"

# --- Row 11 (item 5, new -- this used to be item 2 "URSOFALK") -------------
$ws.Range("A11").Value = 5
Set-TextValue "C11" "URSOFALK 250MG 20 CAPS."
Set-TextValue "H11" "0:0"
Set-TextValue "L11" "1"
Set-TextValue "N11" "122.00"
Set-TextValue "P11" "122.0000"
Set-TextValue "Q11" "1:0"

# --- Row 12 (item 6, new) --------------------------------------------------
$ws.Range("A12").Value = 6
Set-TextValue "C12" "سرنجات 3 سم"
Set-TextValue "H12" "0:0"
Set-TextValue "L12" "0"
Set-TextValue "N12" "2.00"
Set-TextValue "P12" "2.0000"
Set-TextValue "Q12" "1:0"

# --- Row 13 (former row 9): totals row, sum of the six "selling price" ----
$ws.Range("P13").Value = 475.82

# --- Row 14 (former row 10): footer, timestamp refreshed to 10:00 AM ------
Set-TextValue "A14" "Thursday, 17 July, 2025 10:00 AM"
